# Scheduled-runner update: refresh computed profit figures (columns H-N)
# for a batch of recipe rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 875.4286
$ws.Range("I38").Value = 306.22223
$ws.Range("J38").Value = 1900
$ws.Range("K38").Value = 918.66669
$ws.Range("L38").Value = 5700
$ws.Range("M38").Value = -546.66669
$ws.Range("N38").Value = -6444

$ws.Range("H61").Value = 66700100
$ws.Range("I61").Value = 66700100
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 200100300
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -200100128

$ws.Range("H113").Value = 3490.75
$ws.Range("I113").Value = 3527
$ws.Range("J113").Value = 3440
$ws.Range("K113").Value = 3527
$ws.Range("L113").Value = 3440
$ws.Range("M113").Value = -273
$ws.Range("N113").Value = -9948

$ws.Range("H132").Value = 22889.867
$ws.Range("I132").Value = 3260.6
$ws.Range("J132").Value = 91592.3
$ws.Range("K132").Value = 9781.799999999999
$ws.Range("L132").Value = 274776.9
$ws.Range("M132").Value = -7251.799999999999
$ws.Range("N132").Value = -279836.9

$ws.Range("H137").Value = 16419.6
$ws.Range("I137").Value = 14654.667
$ws.Range("J137").Value = 19067
$ws.Range("K137").Value = 43964.001
$ws.Range("L137").Value = 57201
$ws.Range("M137").Value = -41414.001
$ws.Range("N137").Value = -62301

$ws.Range("H138").Value = 2859.375
$ws.Range("I138").Value = 1875
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 5625
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -485
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12250.62
$ws.Range("I32").Value = 11497.357
$ws.Range("J32").Value = 16205.25
$ws.Range("K32").Value = 11497.357
$ws.Range("L32").Value = 16205.25
$ws.Range("M32").Value = -11210.357
$ws.Range("N32").Value = -16779.25

$ws.Range("H61").Value = 2394.64
$ws.Range("I61").Value = 1911.2222
$ws.Range("J61").Value = 3637.7144
$ws.Range("K61").Value = 1911.2222
$ws.Range("L61").Value = 3637.7144
$ws.Range("M61").Value = -1699.2222
$ws.Range("N61").Value = -4061.7144

$ws.Range("H74").Value = 1376.1296
$ws.Range("I74").Value = 1138.2954
$ws.Range("J74").Value = 2422.6
$ws.Range("K74").Value = 1138.2954
$ws.Range("L74").Value = 2422.6
$ws.Range("M74").Value = -264.2954
$ws.Range("N74").Value = -4170.6

$ws.Range("H77").Value = 1376.1296
$ws.Range("I77").Value = 1138.2954
$ws.Range("J77").Value = 2422.6
$ws.Range("K77").Value = 5691.477
$ws.Range("L77").Value = 12113
$ws.Range("M77").Value = -1323.477
$ws.Range("N77").Value = -20849

$ws.Range("H122").Value = 1549.7441
$ws.Range("I122").Value = 1468.2433
$ws.Range("J122").Value = 2052.3333
$ws.Range("K122").Value = 4404.7299
$ws.Range("L122").Value = 6156.999899999999
$ws.Range("M122").Value = -1954.7299
$ws.Range("N122").Value = -11056.9999

$ws.Range("H132").Value = 13160704
$ws.Range("I132").Value = 23811400
$ws.Range("J132").Value = 3961.353
$ws.Range("K132").Value = 71434200
$ws.Range("L132").Value = 11884.059
$ws.Range("M132").Value = -71431670
$ws.Range("N132").Value = -16944.059

$ws.Range("H136").Value = 2394.64
$ws.Range("I136").Value = 1911.2222
$ws.Range("J136").Value = 3637.7144
$ws.Range("K136").Value = 5733.6666
$ws.Range("L136").Value = 10913.1432
$ws.Range("M136").Value = -3183.6666
$ws.Range("N136").Value = -16013.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2785.56
$ws.Range("I134").Value = 2439.1052
$ws.Range("J134").Value = 3882.6667
$ws.Range("K134").Value = 7317.3156
$ws.Range("L134").Value = 11648.0001
$ws.Range("M134").Value = -4782.3156
$ws.Range("N134").Value = -16718.0001

$ws.Range("H137").Value = 55842.082
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 55842.082
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 55842.082
$ws.Range("N137").Value = -66042.08199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8135180
$ws.Range("I31").Value = 2818.3845
$ws.Range("J31").Value = 11910919
$ws.Range("K31").Value = 2818.3845
$ws.Range("L31").Value = 11910919
$ws.Range("M31").Value = -2523.3845
$ws.Range("N31").Value = -11911509

$ws.Range("H34").Value = 8135180
$ws.Range("I34").Value = 2818.3845
$ws.Range("J34").Value = 11910919
$ws.Range("K34").Value = 2818.3845
$ws.Range("L34").Value = 11910919
$ws.Range("M34").Value = -2616.3845
$ws.Range("N34").Value = -11911323

$ws.Range("H58").Value = 14287463
$ws.Range("I58").Value = 981.7222
$ws.Range("J58").Value = 29414324
$ws.Range("K58").Value = 981.7222
$ws.Range("L58").Value = 29414324
$ws.Range("M58").Value = -778.7222
$ws.Range("N58").Value = -29414730

$ws.Range("H99").Value = 1642.4166
$ws.Range("I99").Value = 1500.9
$ws.Range("J99").Value = 2350
$ws.Range("K99").Value = 1500.9
$ws.Range("L99").Value = 2350
$ws.Range("M99").Value = -2.900000000000091
$ws.Range("N99").Value = -5346

$ws.Range("H126").Value = 1642.4166
$ws.Range("I126").Value = 1500.9
$ws.Range("J126").Value = 2350
$ws.Range("K126").Value = 4502.700000000001
$ws.Range("L126").Value = 7050
$ws.Range("M126").Value = -2032.700000000001
$ws.Range("N126").Value = -11990

$ws.Range("H132").Value = 35175.14
$ws.Range("I132").Value = 1614.8125
$ws.Range("J132").Value = 132805.19
$ws.Range("K132").Value = 4844.4375
$ws.Range("L132").Value = 398415.57
$ws.Range("M132").Value = -2314.4375
$ws.Range("N132").Value = -403475.57

$ws.Range("H134").Value = 292888.94
$ws.Range("I134").Value = 1008.2895
$ws.Range("J134").Value = 1402035.4
$ws.Range("K134").Value = 3024.8685
$ws.Range("L134").Value = 4206106.199999999
$ws.Range("M134").Value = -489.8685
$ws.Range("N134").Value = -4211176.199999999

$ws.Range("H136").Value = 14287463
$ws.Range("I136").Value = 981.7222
$ws.Range("J136").Value = 29414324
$ws.Range("K136").Value = 2945.1666
$ws.Range("L136").Value = 88242972
$ws.Range("M136").Value = -395.1666
$ws.Range("N136").Value = -88248072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6019546.5
$ws.Range("I68").Value = 6536812
$ws.Range("J68").Value = 5815046.5
$ws.Range("K68").Value = 19610436
$ws.Range("L68").Value = 17445139.5
$ws.Range("M68").Value = -19609625
$ws.Range("N68").Value = -17446761.5

$ws.Range("H71").Value = 6019546.5
$ws.Range("I71").Value = 6536812
$ws.Range("J71").Value = 5815046.5
$ws.Range("K71").Value = 58831308
$ws.Range("L71").Value = 52335418.5
$ws.Range("M71").Value = -58827252
$ws.Range("N71").Value = -52343530.5

$ws.Range("H92").Value = 2558913
$ws.Range("I92").Value = 11765560
$ws.Range("J92").Value = 1511.1666
$ws.Range("K92").Value = 35296680
$ws.Range("L92").Value = 4533.4998
$ws.Range("M92").Value = -35295432
$ws.Range("N92").Value = -7029.4998

$ws.Range("H131").Value = 2818.492
$ws.Range("I131").Value = 13167.625
$ws.Range("J131").Value = 1313.1637
$ws.Range("K131").Value = 39502.875
$ws.Range("L131").Value = 3939.4911
$ws.Range("M131").Value = -34462.875
$ws.Range("N131").Value = -14019.4911

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1520.1364
$ws.Range("I102").Value = 1044.9286
$ws.Range("J102").Value = 2351.75
$ws.Range("K102").Value = 1044.9286
$ws.Range("L102").Value = 2351.75
$ws.Range("M102").Value = 577.0714
$ws.Range("N102").Value = -5595.75

$ws.Range("H132").Value = 3660.889
$ws.Range("I132").Value = 2529.7
$ws.Range("J132").Value = 5074.875
$ws.Range("K132").Value = 7589.099999999999
$ws.Range("L132").Value = 15224.625
$ws.Range("M132").Value = -5059.099999999999
$ws.Range("N132").Value = -20284.625

$ws.Range("H134").Value = 13969.8
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 13969.8
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 41909.39999999999
$ws.Range("N134").Value = -46979.39999999999

$ws.Range("H140").Value = 36034.75
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 36034.75
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 36034.75
$ws.Range("N140").Value = -46394.75

$ws.Range("H141").Value = 70098.75
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 70098.75
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 70098.75
$ws.Range("N141").Value = -80458.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3636.0293
$ws.Range("I132").Value = 2492.2
$ws.Range("J132").Value = 5270.0713
$ws.Range("K132").Value = 7476.599999999999
$ws.Range("L132").Value = 15810.2139
$ws.Range("M132").Value = -4946.599999999999
$ws.Range("N132").Value = -20870.2139

$ws.Range("H136").Value = 1989.697
$ws.Range("I136").Value = 1244.409
$ws.Range("J136").Value = 3480.2727
$ws.Range("K136").Value = 3733.227
$ws.Range("L136").Value = 10440.8181
$ws.Range("M136").Value = -1183.227
$ws.Range("N136").Value = -15540.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 15000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 15000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15496

$ws.Range("H67").Value = 15000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 15000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16716

$ws.Range("H126").Value = 952.75
$ws.Range("I126").Value = 984.3333
$ws.Range("J126").Value = 912.1429000000001
$ws.Range("K126").Value = 2952.9999
$ws.Range("L126").Value = 2736.4287
$ws.Range("M126").Value = -482.9998999999998
$ws.Range("N126").Value = -7676.4287

$ws.Range("H132").Value = 26318594
$ws.Range("I132").Value = 2208.0908
$ws.Range("J132").Value = 62503624
$ws.Range("K132").Value = 6624.2724
$ws.Range("L132").Value = 187510872
$ws.Range("M132").Value = -4094.2724
$ws.Range("N132").Value = -187515932

$ws.Range("H140").Value = 43038.734
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 43038.734
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 43038.734
$ws.Range("N140").Value = -53398.734
